# Apply the Sprint Backlog edits described in the commit:
#  - Fix "Moduo do sistema" typo -> "Modulo do sistema"
#  - Clarify "Issue" -> "Issue (serviço)"
#  - Fill in the WorkLog ("F") column, which was left blank for most rows:
#      rows 4-19 -> "4h" (matching the value already present in F3)
#      rows 20-22 -> "3h" (new value)
#  - Move the active selection from F3 to B2

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SprintBacklog")

# Header row fixes
$ws.Range("B2").Value = "Modulo do sistema"
$ws.Range("C2").Value = "Issue (serviço)"

# Fill in missing WorkLog values
$ws.Range("F4:F19").Value = "4h"
$ws.Range("F20:F22").Value = "3h"

# Move the selection / active cell to B2
$ws.Activate()
$ws.Range("B2").Select()
